# Apply updated crypto symbol list values (refresh run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'243.64"
$ws.Range('G2').Value = "'9"
$ws.Range('G3').Value = "'9"
$ws.Range('B4').Value = "LEO"
$ws.Range('C4').Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range('D4').Value = "'3.578"
$ws.Range('E4').Value = "3LEOLEO"
$ws.Range('G4').Value = "'9"
$ws.Range('B5').Value = "HuobiToken"
$ws.Range('C5').Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('D5').Value = "'5.300"
$ws.Range('E5').Value = "4HuobiTokenHT"
$ws.Range('G5').Value = "'9"
$ws.Range('B6').Value = "Cronos"
$ws.Range('C6').Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('D6').Value = "'0.05793"
$ws.Range('E6').Value = "5CronosCRO"
$ws.Range('G6').Value = "'9"
$ws.Range('B7').Value = "KuCoinToken"
$ws.Range('C7').Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range('D7').Value = "'6.491"
$ws.Range('E7').Value = "6KuCoinTokenKCS"
$ws.Range('G7').Value = "'9"
$ws.Range('B8').Value = "GateToken"
$ws.Range('C8').Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range('D8').Value = "'3.342"
$ws.Range('E8').Value = "7GateTokenGT"
$ws.Range('G8').Value = "'9"
$ws.Range('B9').Value = "MXToken"
$ws.Range('C9').Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('D9').Value = "'0.8086"
$ws.Range('E9').Value = "8MXTokenMX"
$ws.Range('G9').Value = "'9"
$ws.Range('B10').Value = "FTXToken"
$ws.Range('C10').Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range('D10').Value = "'0.8763"
$ws.Range('E10').Value = "9FTXTokenFTT"
$ws.Range('G10').Value = "'9"
$ws.Range('B11').Value = "One"
$ws.Range('C11').Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range('D11').Value = "'0.01041"
$ws.Range('E11').Value = "10OneONEBestin24h"
$ws.Range('G11').Value = "'9"
$ws.Range('B12').Value = "WazirX"
$ws.Range('C12').Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range('D12').Value = "'0.1388"
$ws.Range('E12').Value = "11WazirXWRX"
$ws.Range('G12').Value = "'9"
$ws.Range('B13').Value = "MandalaExchangeToken"
$ws.Range('C13').Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range('D13').Value = "'0.07276"
$ws.Range('E13').Value = "12MandalaExchangeTokenMDX"
$ws.Range('G13').Value = "'9"
$ws.Range('B14').Value = "LiechtensteinCryptoassetsExchange"
$ws.Range('C14').Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range('D14').Value = "'0.03075"
$ws.Range('E14').Value = "13LiechtensteinCryptoassetsExchangeLCX"
$ws.Range('G14').Value = "'9"
$ws.Range('B15').Value = "BitrueCoin"
$ws.Range('C15').Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range('D15').Value = "'0.03057"
$ws.Range('E15').Value = "14BitrueCoinBTR"
$ws.Range('G15').Value = "'9"
$ws.Range('B16').Value = "BitMartToken"
$ws.Range('C16').Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range('D16').Value = "'0.09325"
$ws.Range('E16').Value = "15BitMartTokenBMX"
$ws.Range('G16').Value = "'9"
$ws.Range('B17').Value = "MCDex"
$ws.Range('C17').Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range('D17').Value = "'3.855"
$ws.Range('E17').Value = "16MCDexMCB"
$ws.Range('G17').Value = "'9"
$ws.Range('B18').Value = "BitForexToken"
$ws.Range('C18').Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range('D18').Value = "'0.001538"
$ws.Range('E18').Value = "17BitForexTokenBF"
$ws.Range('G18').Value = "'9"
$ws.Range('B19').Value = "CoinExToken"
$ws.Range('C19').Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range('D19').Value = "'0.04717"
$ws.Range('E19').Value = "18CoinExTokenCET"
$ws.Range('G19').Value = "'9"
$ws.Range('B20').Value = "TigerCash"
$ws.Range('C20').Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range('D20').Value = "'0.006015"
$ws.Range('E20').Value = "19TigerCashTCH"
$ws.Range('G20').Value = "'9"
$ws.Range('B21').Value = "BitKan"
$ws.Range('C21').Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range('D21').Value = "'0.001269"
$ws.Range('E21').Value = "20BitKanKAN"
$ws.Range('G21').Value = "'9"
$ws.Range('B22').Value = "HotbitToken"
$ws.Range('C22').Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range('D22').Value = "'0.004593"
$ws.Range('E22').Value = "21HotbitTokenHTB"
$ws.Range('G22').Value = "'9"
$ws.Range('B23').Value = "NitroEx"
$ws.Range('C23').Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range('D23').Value = "'0.00008701"
$ws.Range('E23').Value = "22NitroExNTX"
$ws.Range('G23').Value = "'9"
$ws.Range('D24').Value = "'2.141"
$ws.Range('G24').Value = "'9"
$ws.Range('D25').Value = "'0.3210"
$ws.Range('G25').Value = "'9"
$ws.Range('G26').Value = "'9"
$ws.Range('G27').Value = "'9"
$ws.Range('D28').Value = "'0.0002344"
$ws.Range('G28').Value = "'9"
$ws.Range('G29').Value = "'9"
$ws.Range('G30').Value = "'9"
$ws.Range('G31').Value = "'9"
$ws.Range('G32').Value = "'9"
$ws.Range('G33').Value = "'9"
$ws.Range('G34').Value = "'9"
$ws.Range('G35').Value = "'9"
$ws.Range('G36').Value = "'9"
$ws.Range('G37').Value = "'9"
$ws.Range('G38').Value = "'9"
$ws.Range('G39').Value = "'9"
$ws.Range('D40').Value = "'0.03784"
$ws.Range('G40').Value = "'9"
$ws.Range('D41').Value = "'0.006231"
$ws.Range('G41').Value = "'9"
$ws.Range('B42').Value = "CEJI"
$ws.Range('C42').Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range('D42').Value = "'0.004001"
$ws.Range('E42').Value = "41CEJICEJI"
$ws.Range('G42').Value = "'9"
$ws.Range('B43').Value = "BKEXToken"
$ws.Range('C43').Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range('D43').Value = "'0.1052"
$ws.Range('E43').Value = "42BKEXTokenBKK"
$ws.Range('G43').Value = "'9"
$ws.Range('D44').Value = "'0.007869"
$ws.Range('G44').Value = "'9"
$ws.Range('D45').Value = "'0.00005478"
$ws.Range('G45').Value = "'9"
$ws.Range('G46').Value = "'9"
$ws.Range('D47').Value = "'0.5501"
$ws.Range('E47').Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range('G47').Value = "'9"
$ws.Range('D48').Value = "'0.006516"
$ws.Range('E48').Value = "47BOLOBOLO"
$ws.Range('G48').Value = "'9"
$ws.Range('D49').Value = "'0.00002100"
$ws.Range('G49').Value = "'9"
$ws.Range('D50').Value = "'0.0002000"
$ws.Range('G50').Value = "'9"
$ws.Range('G51').Value = "'9"
